# Fix issue with same user having multiple classes.
#
# Melissa Kimball already has a row (row 6) for COMP 474, but she is also
# enrolled in COMP 472 and needs her own row for that class's grades.
# Duplicate her existing row (to keep the same data types/formatting for
# the ID Number column, etc.) and then update just the class/course cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate Melissa Kimball's row (row 6) into the new row 7.
$ws.Range("A6:G6").Copy($ws.Range("A7:G7"))

# Row 7 is for her COMP 472 class (grade/retake grade stay the same: D/F).
$ws.Range("E7").Value = "COMP 472"

# Restore the active cell selection to where the user left off editing.
$ws.Range("F11").Select()
